# Add two new columns, I (I0) and J (IF), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from an existing
# header cell (H1) onto the two new header cells before setting values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns (rows 2-8)
$i0 = @(1, 6, 5, 1, 9, 7, 5)
$if = @(3, 8, 7, 6, 9, 8, 8)

for ($r = 0; $r -lt 7; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}

Write-Output "done"
